$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# "Volume 31   Number  29" -> "Volume 31   Number  30"
$ws.Range("A8").Value = "Volume 31   Number  30"

# "Report Covering the Week  7/15/2024  Through  7/21/2024"
# -> "Report Covering the Week  7/22/2024  Through  7/28/2024"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Crime statistics table updates (rows 14-30, 33) ---
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = '#,##0'
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = -70
$ws.Range("L14").Value = -75
$ws.Range("M14").Value = 200
$ws.Range("M14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N14").Value = -86.363636363636
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = -35
$ws.Range("L15").Value = -23.529411764705
$ws.Range("M15").Value = 0
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -35.714285714285
$ws.Range("F16").Value = 42
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 282
$ws.Range("J16").Value = 228
$ws.Range("K16").Value = 23.684210526315
$ws.Range("L16").Value = 7.633587786259
$ws.Range("M16").Value = 40.298507462686
$ws.Range("N16").Value = -56.880733944954
$ws.Range("C17").Value = 21
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 59
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = -1.666666666666
$ws.Range("I17").Value = 430
$ws.Range("J17").Value = 383
$ws.Range("K17").Value = 12.271540469973
$ws.Range("L17").Value = 15.281501340482
$ws.Range("M17").Value = 110.78431372549
$ws.Range("N17").Value = 0.702576112412
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 14.285714285714
$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 21.739130434782
$ws.Range("I18").Value = 181
$ws.Range("J18").Value = 133
$ws.Range("K18").Value = 36.090225563909
$ws.Range("L18").Value = 15.286624203821
$ws.Range("M18").Value = 37.121212121212
$ws.Range("N18").Value = -76.854219948849
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 17.142857142857
$ws.Range("I19").Value = 343
$ws.Range("J19").Value = 275
$ws.Range("K19").Value = 24.727272727272
$ws.Range("L19").Value = 29.433962264150
$ws.Range("M19").Value = 86.413043478260
$ws.Range("N19").Value = 36.653386454183
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -11.111111111111
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = -18.75
$ws.Range("I20").Value = 170
$ws.Range("J20").Value = 198
$ws.Range("K20").Value = -14.141414141414
$ws.Range("L20").Value = 4.938271604938
$ws.Range("M20").Value = 150
$ws.Range("N20").Value = -48.328267477203
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 54
$ws.Range("E21").Value = 3.703703703703
$ws.Range("F21").Value = 197
$ws.Range("G21").Value = 195
$ws.Range("H21").Value = 1.025641025641
$ws.Range("I21").Value = 1422
$ws.Range("J21").Value = 1247
$ws.Range("K21").Value = 14.033680834001
$ws.Range("L21").Value = 13.942307692307
$ws.Range("M21").Value = 77.085927770859
$ws.Range("N21").Value = -43.074459567654
$ws.Range("M22").Value = -33.333333333333
$ws.Range("C23").Value = 2
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = 33.333333333333
$ws.Range("L23").Value = 11.111111111111
$ws.Range("M23").Value = 122.222222222222
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 4.545454545454
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = -9.876543209876
$ws.Range("I24").Value = 579
$ws.Range("J24").Value = 606
$ws.Range("K24").Value = -4.455445544554
$ws.Range("L24").Value = -16.810344827586
$ws.Range("M24").Value = 3.208556149732
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -11.111111111111
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = -35.135135135135
$ws.Range("I25").Value = 159
$ws.Range("J25").Value = 215
$ws.Range("K25").Value = -26.046511627907
$ws.Range("L25").Value = -49.363057324840
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 28
$ws.Range("E26").Value = -46.428571428571
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 88
$ws.Range("H26").Value = -5.681818181818
$ws.Range("I26").Value = 568
$ws.Range("J26").Value = 510
$ws.Range("K26").Value = 11.372549019607
$ws.Range("L26").Value = 8.396946564885
$ws.Range("M26").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -24.242424242424
$ws.Range("L27").Value = -28.571428571428
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -36.363636363636
$ws.Range("I28").Value = 42
$ws.Range("J28").Value = 46
$ws.Range("K28").Value = -8.695652173913
$ws.Range("L28").Value = -12.5
$ws.Range("C29").Value = 3
$ws.Range("F29").Value = 6
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 24
$ws.Range("K29").Value = 14.285714285714
$ws.Range("L29").Value = -36.842105263157
$ws.Range("M29").Value = 4.347826086956
$ws.Range("N29").Value = -69.230769230769
$ws.Range("C30").Value = 2
$ws.Range("F30").Value = 5
$ws.Range("H30").Value = 150
$ws.Range("I30").Value = 21
$ws.Range("K30").Value = 50
$ws.Range("L30").Value = -34.375
$ws.Range("M30").Value = 5
$ws.Range("N30").Value = -68.656716417910
$ws.Range("L33").Value = -50
